$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

# Update the "Closure of schools/universities" note for Italy (row 30, column E)
$ws.Range("E30").Value = "Closure of schools and universities from March 4 until April 3, extended to at least April 14."

# Update the "Cancellation of public events / Closure of public places" note for Italy (row 30, column F)
$ws.Range("F30").Value = "- Bars and restaurants along with many other retail trade activities (e.g. shopping centres; indoor and outdoor markets) closed from March 10 until at least 14 April, and all sporting competitions suspended over the same period along with other public gatherings.`n- All but prescribed essential production activities suspended from March 23, with the list of permitted activities further limited from March 26.`n- On March 30, closures extended from April 3 to 30 April for sports, bars and similar activities. "

# Bump the "Updated on" date column (B) by one day for every populated row (5-74)
for ($r = 5; $r -le 74; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v + 1
    }
}
